# Updated cryptos list with GitHub Actions
# Refreshes the Price (column D) and Volume(1h) (column E) figures for the
# coin rows, and reorders the TheSandbox / Quant rows (41 <-> 42) to match
# the new ranking, each keeping the rest of the row layout intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param(
        $ws,
        [string]$cellRef,
        [string]$val
    )
    $rng = $ws.Range($cellRef)
    # Force the cell to stay plain text even when the value looks numeric
    # (e.g. "0.9995"), mirroring the inline-string cells in the source file.
    $rng.NumberFormat = "@"
    $rng.Value = $val
    # Drop back to the default style so no stray formatting/quote-prefix
    # is left behind on the cell.
    $rng.Style = "Normal"
}

Set-CellText $ws 'D2' '30.533.89'
Set-CellText $ws 'E2' '  -0.25%  '
Set-CellText $ws 'D3' '1.911.68'
Set-CellText $ws 'E3' '  -0.59%  '
Set-CellText $ws 'D4' '0.9995'
Set-CellText $ws 'E4' '  -0.06%  '
Set-CellText $ws 'D5' '244.23'
Set-CellText $ws 'E5' '  -1.25%  '
Set-CellText $ws 'D6' '0.9997'
Set-CellText $ws 'E6' '  -0.03%  '
Set-CellText $ws 'D7' '0.4838'
Set-CellText $ws 'E7' '  +1.94%  '
Set-CellText $ws 'D8' '0.2892'
Set-CellText $ws 'E8' '  +0.08%  '
Set-CellText $ws 'D9' '0.06816'
Set-CellText $ws 'E9' '  -0.29%  '
Set-CellText $ws 'D10' '111.24'
Set-CellText $ws 'E10' '  +5.83%  '
Set-CellText $ws 'D11' '19.28'
Set-CellText $ws 'E11' '  +4.93%  '
Set-CellText $ws 'D12' '1.918.10'
Set-CellText $ws 'E12' '  -0.26%  '
Set-CellText $ws 'D13' '0.07569'
Set-CellText $ws 'E13' '  -1.70%  '
Set-CellText $ws 'D14' '5.380'
Set-CellText $ws 'E14' '  +1.19%  '
Set-CellText $ws 'D15' '0.6704'
Set-CellText $ws 'E15' '  +0.40%  '
Set-CellText $ws 'D16' '295.30'
Set-CellText $ws 'E16' '  +1.07%  '
Set-CellText $ws 'D17' '30.530.29'
Set-CellText $ws 'D18' '13.04'
Set-CellText $ws 'E18' '  +0.66%  '
Set-CellText $ws 'D19' '0.9999'
Set-CellText $ws 'D20' '0.000007586'
Set-CellText $ws 'E20' '  -0.49%  '
Set-CellText $ws 'D21' '5.523'
Set-CellText $ws 'E21' '  -1.32%  '
Set-CellText $ws 'D22' '2.159.62'
Set-CellText $ws 'E22' '  -0.55%  '
Set-CellText $ws 'D23' '0.9998'
Set-CellText $ws 'E23' '  -0.07%  '
Set-CellText $ws 'D24' '6.446'
Set-CellText $ws 'E24' '  +0.24%  '
Set-CellText $ws 'D25' '9.464'
Set-CellText $ws 'E25' '  +0.39%  '
Set-CellText $ws 'D26' '165.63'
Set-CellText $ws 'E26' '  -1.31%  '
Set-CellText $ws 'D27' '20.32'
Set-CellText $ws 'E27' '  -3.57%  '
Set-CellText $ws 'D28' '2.081'
Set-CellText $ws 'E28' '  -1.81%  '
Set-CellText $ws 'D29' '0.1064'
Set-CellText $ws 'E29' '  -0.82%  '
Set-CellText $ws 'D30' '1.436'
Set-CellText $ws 'E30' '  +2.94%  '
Set-CellText $ws 'D31' '4.145'
Set-CellText $ws 'E31' '  -0.87%  '
Set-CellText $ws 'D32' '4.053'
Set-CellText $ws 'E32' '  -0.41%  '
Set-CellText $ws 'D33' '0.04971'
Set-CellText $ws 'E33' '  -1.50%  '
Set-CellText $ws 'D34' '0.7371'
Set-CellText $ws 'E34' '  -0.15%  '
Set-CellText $ws 'D35' '1.135'
Set-CellText $ws 'E35' '  -0.80%  '
Set-CellText $ws 'D36' '0.9996'
Set-CellText $ws 'E36' '  +0.05%  '
Set-CellText $ws 'D37' '0.02041'
Set-CellText $ws 'E37' '  -1.77%  '
Set-CellText $ws 'D38' '2.716'
Set-CellText $ws 'E38' '  -0.80%  '
Set-CellText $ws 'D39' '2.684'
Set-CellText $ws 'E39' '  -0.31%  '
Set-CellText $ws 'E40' '  -1.84%  '
Set-CellText $ws 'B41' 'TheSandbox'
Set-CellText $ws 'C41' 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-CellText $ws 'D41' '0.4469'
Set-CellText $ws 'E41' '  +1.75%  '
Set-CellText $ws 'B42' 'Quant'
Set-CellText $ws 'C42' 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-CellText $ws 'D42' '109.29'
Set-CellText $ws 'E42' '  -1.75%  '
Set-CellText $ws 'D43' '0.8670'
Set-CellText $ws 'E43' '  -0.93%  '
Set-CellText $ws 'D44' '5.776'
Set-CellText $ws 'E44' '  -1.90%  '
Set-CellText $ws 'D45' '0.9995'
Set-CellText $ws 'E45' '  -0.04%  '
Set-CellText $ws 'D46' '69.42'
Set-CellText $ws 'E46' '  +2.02%  '
Set-CellText $ws 'D47' '7.208'
Set-CellText $ws 'E47' '  -0.86%  '
Set-CellText $ws 'D48' '48.21'
Set-CellText $ws 'E48' '  -0.32%  '
Set-CellText $ws 'D49' '9.225'
Set-CellText $ws 'E49' '  -1.59%  '
Set-CellText $ws 'E50' '  -1.33%  '
Set-CellText $ws 'E51' '  -0.24%  '
